$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 17:50"

# Full country dataset in final (post-update) sort order: Pais, Casos totales,
# Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$countries = @(
    @("Estados Unidos", 145131, 1640, 4574, 137949, 2972, 26, 2608),
    @("Italia", 97689, 0, 13030, 73880, 3906, 0, 10779),
    @("España", 85195, 5085, 16780, 61075, 5231, 537, 7340),
    @("China", 81470, 31, 75700, 2466, 633, 4, 3304),
    @("Alemania", 63929, 1494, 9211, 54158, 1979, 19, 560),
    @("Iran", 41495, 3186, 13911, 24827, 3511, 117, 2757),
    @("Francia", 40174, 0, 7202, 30366, 4632, 0, 2606),
    @("Reino Unido", 22141, 2619, 135, 20598, 163, 180, 1408),
    @("Suiza", 15526, 697, 1823, 13370, 301, 33, 333),
    @("Belgica", 11899, 1063, 1527, 9859, 927, 82, 513),
    @("Paises Bajos", 11750, 884, 250, 10636, 1053, 93, 864),
    @("Corea del Sur", 9661, 78, 5228, 4275, 59, 6, 158),
    @("Austria", 9454, 666, 636, 8710, 193, 22, 108),
    @("Turquia", 9217, 0, 105, 8981, 568, 0, 131),
    @("Canada", 6671, 351, 1014, 5590, 120, 2, 67),
    @("Portugal", 6408, 446, 43, 6225, 164, 21, 140),
    @("Noruega", 4436, 152, 12, 4392, 97, 6, 32),
    @("Brasil", 4362, 106, 120, 4101, 296, 5, 141),
    @("Israel", 4347, 100, 134, 4197, 80, 1, 16),
    @("Australia", 4245, 82, 244, 3983, 28, 1, 18),
    @("Suecia", 4028, 328, 16, 3866, 306, 36, 146),
    @("Chequia", 2896, 79, 11, 2868, 52, 1, 17),
    @("Malasia", 2626, 156, 479, 2110, 94, 2, 37),
    @("Irlanda", 2615, 0, 5, 2564, 59, 0, 46),
    @("Dinamarca", 2577, 182, 1, 2499, 113, 5, 77),
    @("Chile", 2449, 310, 156, 2285, 14, 1, 8),
    @("Luxemburgo", 1988, 38, 40, 1926, 25, 1, 22),
    @("Polonia", 1984, 122, 7, 1951, 3, 4, 26),
    @("Rumania", 1952, 137, 206, 1700, 31, 3, 46),
    @("Ecuador", 1924, 0, 3, 1863, 58, 0, 58),
    @("Japon", 1866, 0, 424, 1388, 56, 0, 54),
    @("Rusia", 1836, 302, 66, 1761, 8, 1, 9),
    @("Pakistan", 1650, 53, 29, 1601, 11, 6, 20),
    @("Filipinas", 1546, 128, 42, 1426, 1, 7, 78),
    @("Tailandia", 1524, 136, 229, 1286, 11, 2, 9),
    @("Arabia Saudita", 1453, 154, 115, 1330, 12, 0, 8),
    @("Indonesia", 1414, 129, 75, 1217, 0, 8, 122),
    @("Finlandia", 1352, 112, 10, 1329, 49, 2, 13),
    @("Sudafrica", 1280, 0, 31, 1247, 7, 0, 2),
    @("Grecia", 1212, 56, 52, 1117, 72, 4, 43),
    @("Islandia", 1086, 66, 157, 927, 25, 0, 2),
    @("India", 1071, 47, 100, 942, 0, 2, 29),
    @("Mexico", 993, 145, 35, 938, 1, 4, 20),
    @("Panama", 989, 0, 4, 961, 36, 0, 24),
    @("Republica Dominicana", 901, 42, 4, 855, 0, 3, 42),
    @("Singapur", 879, 35, 228, 648, 19, 0, 3),
    @("Peru", 852, 0, 16, 818, 40, 0, 18),
    @("Argentina", 820, 0, 228, 570, 0, 2, 22),
    @("Croacia", 790, 77, 67, 717, 27, 0, 6),
    @("Serbia", 785, 44, 42, 727, 62, 3, 16),
    @("Eslovenia", 756, 26, 10, 735, 28, 0, 11),
    @("Estonia", 715, 36, 20, 692, 10, 0, 3),
    @("Crucero", 712, 0, 603, 99, 15, 0, 10),
    @("Colombia", 702, 0, 10, 682, 29, 0, 10),
    @("Hong Kong", 642, 0, 118, 520, 5, 0, 4),
    @("Catar", 634, 0, 48, 585, 6, 0, 1),
    @("Emiratos Arabes Unidos", 611, 41, 61, 545, 2, 2, 5),
    @("Egipto", 609, 0, 132, 437, 0, 0, 40),
    @("Nueva Zelanda", 589, 75, 63, 525, 2, 0, 1),
    @("Irak", 547, 0, 143, 362, 0, 0, 42),
    @("Marruecos", 516, 37, 14, 473, 1, 3, 29),
    @("Barein", 515, 16, 279, 232, 2, 0, 4),
    @("Argelia", 511, 0, 31, 449, 0, 0, 31),
    @("Lituania", 484, 24, 1, 476, 5, 0, 7),
    @("Armenia", 482, 58, 30, 449, 15, 0, 3),
    @("Ucrania", 480, 5, 6, 463, 0, 1, 11),
    @("Hungria", 447, 39, 34, 398, 6, 2, 15),
    @("Libano", 446, 8, 35, 400, 3, 1, 11),
    @("Letonia", 376, 29, 1, 375, 3, 0, 0),
    @("Bulgaria", 359, 13, 17, 334, 13, 0, 8),
    @("Bosnia y Herzegovina", 358, 35, 17, 332, 1, 3, 9),
    @("Eslovaquia", 336, 22, 7, 329, 1, 0, 0),
    @("Principado de Andorra", 334, 0, 6, 322, 10, 0, 6),
    @("Costa Rica", 314, 0, 3, 309, 6, 0, 2),
    @("Tunez", 312, 0, 2, 302, 10, 0, 8),
    @("Uruguay", 310, 6, 0, 309, 9, 0, 1),
    @("Taiwan", 306, 8, 39, 262, 0, 2, 5),
    @("Kazajistan", 302, 18, 21, 280, 0, 0, 1),
    @("Moldavia", 298, 35, 15, 281, 44, 0, 2),
    @("Republica de Macedonia", 285, 26, 12, 266, 1, 1, 7),
    @("Azerbaiyan", 273, 64, 26, 243, 23, 0, 4),
    @("Kuwait", 266, 11, 72, 194, 13, 0, 0),
    @("Jordania", 259, 0, 18, 237, 3, 1, 4),
    @("Republica de Chipre", 230, 16, 15, 208, 3, 1, 7),
    @("San Marino", 230, 6, 13, 192, 16, 3, 25),
    @("Albania", 223, 11, 44, 168, 7, 1, 11),
    @("Burkina Faso", 222, 0, 23, 187, 0, 0, 12),
    @("Vietnam", 203, 9, 55, 148, 3, 0, 0),
    @("Reunion", 183, 0, 1, 182, 0, 0, 0),
    @("Oman", 179, 12, 29, 150, 3, 0, 0),
    @("Cuba", 170, 31, 4, 162, 2, 1, 4),
    @("Islas Feroe", 168, 9, 70, 98, 1, 0, 0),
    @("Costa de Marfil", 165, 0, 4, 160, 0, 0, 1),
    @("Senegal", 162, 20, 27, 135, 0, 0, 0),
    @("Malta", 156, 5, 2, 154, 4, 0, 0),
    @("Ghana", 152, 0, 2, 145, 1, 0, 5),
    @("Bielorrusia", 152, 58, 32, 120, 2, 0, 0),
    @("Uzbekistan", 149, 5, 7, 140, 8, 0, 2),
    @("Afganistan", 145, 25, 2, 139, 0, 0, 4),
    @("Honduras", 139, 29, 3, 133, 4, 0, 3),
    @("Camerun", 139, 0, 5, 128, 0, 0, 6),
    @("Venezuela", 129, 10, 39, 87, 6, 0, 3),
    @("Mauricio", 128, 21, 0, 125, 1, 0, 3),
    @("Brunei", 127, 1, 38, 88, 3, 0, 1),
    @("Sri Lanka", 122, 5, 15, 105, 5, 1, 2),
    @("Estado de Palestina", 115, 6, 18, 96, 0, 0, 1),
    @("Nigeria", 111, 0, 3, 107, 0, 0, 1),
    @("Camboya", 107, 4, 21, 86, 1, 0, 0),
    @("Guadalupe", 106, 0, 17, 85, 10, 0, 4),
    @("Georgia", 100, 9, 18, 82, 6, 0, 0),
    @("Bolivia", 97, 16, 0, 93, 3, 3, 4),
    @("Kirguistan", 94, 10, 3, 91, 0, 0, 0),
    @("Martinica", 93, 0, 0, 92, 12, 0, 1),
    @("Montenegro", 91, 6, 0, 90, 1, 0, 1),
    @("Trinidad yTobago", 82, 4, 1, 78, 0, 0, 3),
    @("Mayotte", 82, 19, 10, 72, 3, 0, 0),
    @("Consejo Danes para los Refugiados", 81, 0, 2, 71, 0, 0, 8),
    @("Ruanda", 70, 0, 0, 70, 0, 0, 0),
    @("Gibraltar", 65, 0, 14, 51, 0, 0, 0),
    @("Paraguay", 64, 5, 1, 60, 3, 0, 3),
    @("Liechtenstein", 62, 6, 0, 62, 0, 0, 0),
    @("Aruba", 50, 0, 1, 49, 0, 0, 0),
    @("Kenia", 50, 8, 1, 48, 2, 0, 1),
    @("Isla de Man", 49, 7, 0, 49, 0, 0, 0),
    @("Banglades", 49, 1, 19, 25, 1, 0, 5),
    @("Monaco", 46, 0, 1, 44, 0, 0, 1),
    @("Guayana Francesa", 43, 0, 6, 37, 0, 0, 0),
    @("Madagascar", 39, 0, 0, 39, 0, 0, 0),
    @("Puerto Rico", 39, 0, 1, 36, 0, 0, 2),
    @("Macao", 38, 1, 10, 28, 0, 0, 0),
    @("Jamaica", 36, 4, 2, 33, 0, 0, 1),
    @("Guatemala", 36, 2, 10, 25, 1, 0, 1),
    @("Polinesia Francesa", 35, 5, 0, 35, 2, 0, 0),
    @("Zambia", 35, 6, 0, 35, 0, 0, 0),
    @("Uganda", 33, 0, 0, 33, 0, 0, 0),
    @("Barbados", 33, 0, 0, 33, 0, 0, 0),
    @("Guam", 32, 0, 0, 31, 0, 0, 1),
    @("El Salvador", 30, 6, 0, 30, 5, 0, 0),
    @("Togo", 30, 5, 1, 28, 0, 0, 1),
    @("Mali", 25, 7, 0, 23, 0, 1, 2),
    @("Etiopia", 23, 2, 4, 19, 1, 0, 0),
    @("Guinea", 22, 6, 0, 22, 0, 0, 0),
    @("Bermudas", 22, 0, 2, 20, 0, 0, 0),
    @("Niger", 22, 4, 0, 19, 0, 2, 3),
    @("Congo", 19, 0, 0, 19, 0, 0, 0),
    @("Tanzania", 19, 5, 1, 18, 0, 0, 0),
    @("Republica de Yibuti", 18, 0, 0, 18, 0, 0, 0),
    @("Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0),
    @("Maldivas", 17, 0, 13, 4, 0, 0, 0),
    @("Nueva Caledonia", 15, 0, 0, 15, 0, 0, 0),
    @("Haiti", 15, 0, 1, 14, 0, 0, 0),
    @("San Martin (Parte Francesa)", 15, 4, 2, 12, 0, 1, 1),
    @("Birmania", 14, 4, 0, 14, 0, 0, 0),
    @("Bahamas", 14, 0, 1, 13, 0, 0, 0),
    @("Guinea Ecuatorial", 12, 0, 0, 12, 0, 0, 0),
    @("Eritrea", 12, 0, 0, 12, 0, 0, 0),
    @("Islas Caimanes", 12, 0, 0, 11, 0, 0, 1),
    @("Mongolia", 12, 0, 2, 10, 0, 0, 0),
    @("Dominica", 11, 0, 0, 11, 0, 0, 0),
    @("Namibia", 11, 0, 2, 9, 0, 0, 0),
    @("Curazao", 11, 3, 2, 8, 0, 0, 1),
    @("Groenlandia", 10, 0, 2, 8, 0, 0, 0),
    @("Granada", 9, 0, 0, 9, 0, 0, 0),
    @("Suazilandia", 9, 0, 0, 9, 0, 0, 0),
    @("Santa Lucia", 9, 0, 1, 8, 0, 0, 0),
    @("Siria", 9, 0, 0, 8, 0, 0, 1),
    @("Laos", 8, 0, 0, 8, 0, 0, 0),
    @("Seychelles", 8, 0, 0, 8, 0, 0, 0),
    @("Surinam", 8, 0, 0, 8, 0, 0, 0),
    @("Mozambique", 8, 0, 0, 8, 0, 0, 0),
    @("Libia", 8, 0, 0, 8, 0, 0, 0),
    @("Guinea-Bisau", 8, 6, 0, 8, 0, 0, 0),
    @("Guyana", 8, 0, 0, 7, 0, 0, 1),
    @("Antigua y Barbuda", 7, 0, 0, 7, 0, 0, 0),
    @("San Cristobal y Nieves", 7, 5, 0, 7, 0, 0, 0),
    @("Zimbabue", 7, 0, 0, 6, 0, 0, 1),
    @("Gabon", 7, 0, 0, 6, 0, 0, 1),
    @("Angola", 7, 0, 0, 5, 0, 0, 2),
    @("Benin", 6, 0, 0, 6, 0, 0, 0),
    @("Santa Sede", 6, 0, 0, 6, 0, 0, 0),
    @("San Martin (Parte Holandesa)", 6, 0, 0, 6, 0, 0, 0),
    @("San Bartolome", 6, 1, 1, 5, 0, 0, 0),
    @("Cabo Verde", 6, 0, 0, 5, 0, 0, 1),
    @("Sudan", 6, 0, 0, 4, 0, 1, 2),
    @("Republica del Chad", 5, 2, 0, 5, 0, 0, 0),
    @("Fiyi", 5, 0, 0, 5, 0, 0, 0),
    @("Islas Turcas y Caicos", 5, 1, 0, 5, 0, 0, 0),
    @("Montserrat", 5, 0, 0, 5, 0, 0, 0),
    @("Nepal", 5, 0, 1, 4, 0, 0, 0),
    @("Mauritania", 5, 0, 2, 3, 0, 0, 0),
    @("Butan", 4, 0, 0, 4, 0, 0, 0),
    @("Nicaragua", 4, 0, 0, 3, 0, 0, 1),
    @("Gambia", 4, 0, 0, 3, 0, 0, 1),
    @("Somalia", 3, 0, 0, 3, 0, 0, 0),
    @("Liberia", 3, 0, 0, 3, 0, 0, 0),
    @("Republica de Africa Central", 3, 0, 0, 3, 0, 0, 0),
    @("Anguila", 2, 0, 0, 2, 0, 0, 0),
    @("Islas Virgenes Britanicas", 2, 0, 0, 2, 0, 0, 0),
    @("Belice", 2, 0, 0, 2, 0, 0, 0),
    @("Timor Oriental", 1, 0, 0, 1, 0, 0, 0),
    @("Papua Nueva Guinea", 1, 0, 0, 1, 0, 0, 0),
    @("San Vicente y las Granadinas", 1, 0, 1, 0, 0, 0, 0)
)

$row = 4
foreach ($c in $countries) {
    $ws.Cells.Item($row, 1).Value = $c[0]
    $ws.Cells.Item($row, 2).Value = $c[1]
    $ws.Cells.Item($row, 3).Value = $c[2]
    $ws.Cells.Item($row, 4).Value = $c[3]
    $ws.Cells.Item($row, 5).Value = $c[4]
    $ws.Cells.Item($row, 6).Value = $c[5]
    $ws.Cells.Item($row, 7).Value = $c[6]
    $ws.Cells.Item($row, 8).Value = $c[7]
    $row++
}
